$wb = $excel.ActiveWorkbook

# --- Sheet 2 (tab) rename: "Include from Evaluation Reaso" -> "Include #0" ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Include #0"

# --- Sheet 1 ("Metadata") edits ---
$ws1 = $wb.Worksheets.Item(1)

# URL (pythia -> cicada)
$ws1.Range("B2").Value = "http://fhirfli.dev/fhir/ig/cicada/ValueSet/eval-reason"

# Date value update
$ws1.Range("B8").Value = "2026-02-11T14:37:07-05:00"

# Insert a new row above the old "Description" row (row 11) for "Jurisdiction",
# keeping the same formatting as the surrounding rows (copy format from row 10).
$ws1.Rows.Item(11).Insert()
$ws1.Range("A10:B10").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122, -4163)
$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = ""

# --- Sheet 2 (content) edits: System URI (pythia -> cicada) ---
$ws2.Range("B18").Value = "http://fhirfli.dev/fhir/ig/cicada/CodeSystem/EvalReason"
